$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 (Item 4) with the new connector part information
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "J5 J6"
$ws.Range("D5").Value = "JST Sales America Inc."
$ws.Range("E5").Value = "B2B-XH-AM(LF)(SN)"
$ws.Range("F5").Value = "1x2pin 100mil"
$ws.Range("G5").Value = "B2B-XH-AM(LF)(SN)"
$ws.Range("H5").Value = "250V"
$ws.Range("I5").Value = "CONN HEADER VERT 2POS 2.5MM"

# Reflect the new (longer) content by widening columns C:G to best-fit values
$ws.Range("C1").ColumnWidth = 13.5
$ws.Range("D1").ColumnWidth = 19.5
$ws.Range("E1").ColumnWidth = 17.7
$ws.Range("F1").ColumnWidth = 12.5
$ws.Range("G1").ColumnWidth = 17.7

$wb.Save()
